$d = $word.ActiveDocument

# Locate the placeholder paragraph ("Tady bude videonávod, až bude hotový :)")
# inside the "Videonávod" section and replace it with the real copy plus a
# new paragraph that will hold the YouTube link.
$target = $d.Content.Find
$target.ClearFormatting()
$found = $d.Content.Find.Execute("Tady bude videonávod, až bude hotový :)")
if (-not $found) {
    throw "placeholder paragraph not found"
}

$placeholder = $d.Content.Duplicate
$placeholder.Start = $d.Content.Find.Parent.Start
$placeholder.End = $d.Content.Find.Parent.End

# Re-locate the paragraph range reliably via the Paragraphs collection so we
# operate on the whole paragraph (incl. its paragraph mark).
$videoPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i).Range
    if ($cand.Text -like "Tady bude videonávod*") {
        $videoPara = $cand
        break
    }
}
if ($null -eq $videoPara) {
    throw "could not re-locate placeholder paragraph"
}

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Videonávod ukazuje celý proces, který se týká správy projektů a akcí od zápisu projektu po jeho uzavření.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Video nemá zvuk.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$videoPara.InsertXML($xmlFrag)

# The InsertXML above produced two new paragraphs in place of the old one;
# the second (currently empty) paragraph is where the YouTube hyperlink goes.
$linkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Videonávod ukazuje*") {
        $linkPara = $d.Paragraphs($i + 1)
        break
    }
}
if ($null -eq $linkPara) {
    throw "could not locate the new empty BodyText paragraph"
}

$insStart = $linkPara.Range.Start
$insPoint = $d.Range($insStart, $insStart)
$url = "https://youtu.be/eJ_H6Tonj9Y?si=zyCa4qXBjETU_jkJ"
$d.Hyperlinks.Add($insPoint, $url, $null, $null, $url) | Out-Null
